$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A; this shifts the existing A:D columns
# (and their widths/data) to B:E and the data in row1-3 accordingly.
$ws.Columns.Item(1).Insert()

# New header cell for the inserted "Identificador" column, bold black font.
$ws.Range("A1").Value = "Identificador"
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").Font.Color = 0

# Row identifiers for the two data rows.
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2

# Give the new column an explicit width (character width 12 once stored).
$ws.Columns.Item(1).ColumnWidth = 11.166666666666666

# Update selection to match the edited workbook state.
$ws.Range("A3").Select()
